# Update the team-specific transition matrix on Sheet1 with the
# recomputed probabilities after adding more games / speeding up the
# simulation logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Af0)
$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.5384615384615384
$ws.Range("J2").Value = 0.03846153846153846
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.07692307692307693

# Row 3 (Af1)
$ws.Range("J3").Value = 0.1428571428571428
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.1428571428571428

# Row 4 (Af2)
$ws.Range("P4").Value = 1

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.1333333333333333
$ws.Range("O6").Value = 0.06666666666666667
$ws.Range("Q6").Value = 0.3333333333333333
$ws.Range("R6").Value = 0.2
$ws.Range("S6").Value = 0.2

# Row 7 (Ai1)
$ws.Range("D7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.1111111111111111
$ws.Range("O7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("R7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.4444444444444444

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.16
$ws.Range("F8").Value = 0.12
$ws.Range("J8").Value = 0.12
$ws.Range("O8").Value = 0.04
$ws.Range("Q8").Value = 0.16
$ws.Range("R8").Value = 0.08
$ws.Range("S8").Value = 0.32

# Row 9 (Ai3)
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.1818181818181818
$ws.Range("Q9").Value = 0.2727272727272727
$ws.Range("R9").Value = 0.2727272727272727
$ws.Range("S9").Value = 0.1818181818181818

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.1311475409836066
$ws.Range("D10").Value = 0.00819672131147541
$ws.Range("O10").Value = 0.00819672131147541
$ws.Range("Q10").Value = 0.180327868852459
$ws.Range("R10").Value = 0.139344262295082
$ws.Range("S10").Value = 0.3360655737704918

# Row 11 (Bf0)
$ws.Range("J11").Value = 0.2666666666666667
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5333333333333333

# Row 13 (Bf2)
$ws.Range("G13").Value = 1

# Row 15 (Bi0)
$ws.Range("H15").Value = 0.09090909090909091
$ws.Range("I15").Value = 0.09090909090909091
$ws.Range("J15").Value = 0.5454545454545454
$ws.Range("O15").Value = 0.09090909090909091
$ws.Range("S15").Value = 0.1818181818181818

# Row 16 (Bi1)
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.75

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.05555555555555555
$ws.Range("H17").Value = 0.02777777777777778
$ws.Range("I17").Value = 0.08333333333333333
$ws.Range("J17").Value = 0.5555555555555556
$ws.Range("K17").Value = 0.1111111111111111
$ws.Range("M17").Value = 0.02777777777777778
$ws.Range("O17").Value = 0.05555555555555555
$ws.Range("S17").Value = 0.08333333333333333

# Row 18 (Bi3)
$ws.Range("H18").Value = 0.1538461538461539
$ws.Range("I18").Value = 0.1153846153846154
$ws.Range("O18").Value = 0.03846153846153846

# Row 19 (Br0)
$ws.Range("H19").Value = 0.2207792207792208
$ws.Range("I19").Value = 0.02597402597402598
$ws.Range("J19").Value = 0.5194805194805194
$ws.Range("K19").Value = 0.07792207792207792
$ws.Range("M19").Value = 0.01298701298701299
$ws.Range("O19").Value = 0.03896103896103896
$ws.Range("S19").Value = 0.1038961038961039
